$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove all existing hyperlinks (cells keep their text); we will
#    re-create the 6 that should remain, in the right order, further
#    down. (Range.Hyperlinks.Delete() clears the whole sheet's set.)
# ------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2. Wipe the third "valid email / empty password" table (G1:H4, I5)
#    and the old rows 5-7 - they get replaced by the new layout.
# ------------------------------------------------------------------
$ws.Range("G1:H4").Clear()
$ws.Range("I5").Clear()
$ws.Range("A5:I7").Clear()

# ------------------------------------------------------------------
# 3. Update the values that changed in place (rows 1-4).
# ------------------------------------------------------------------
$ws.Range("A4").Value = "abc@gmailcom"

$ws.Range("D2").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("D4").Value = " "

$ws.Range("E2").Value = "halo1"
$ws.Range("E3").Value = "batangcoklat"
$ws.Range("E4").Value = "halo3"

# ------------------------------------------------------------------
# 4. G3 becomes a plain bordered blank cell (same look as I3), and a
#    new plain bordered blank cell appears at B5. (PasteSpecial with
#    "formats only" never touches the value, so no value ever lands
#    in either cell - no need for an extra ClearContents afterwards,
#    which was observed to drop the cell/row entirely.)
# ------------------------------------------------------------------
$ws.Range("I3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 5. New header row 7: "email" / "invalid password", styled like the
#    other header rows. Copy column-by-column (single-cell copies)
#    to avoid the paste-special engine minting a duplicate near-
#    identical cellXf when a multi-cell range is copied as a block.
# ------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A7").Value = "email"
$ws.Range("B7").Value = "invalid password"

# ------------------------------------------------------------------
# 6. New data rows 8-10 (email / password pairs), styled like the
#    existing A3 / B3 data row (hyperlink font + plain bordered
#    font), again pasted one column at a time.
# ------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A8").Value = "titovalvala@gmail.com"
$ws.Range("B8").Value = "password1"

$ws.Range("A9").Value = "valiantartwear@gmail.com"
$ws.Range("B9").Value = "1213651%%"

$ws.Range("A10").Value = "titovalvala@gmail.com"
$ws.Range("B10").Value = "HAHAHUHU"

# ------------------------------------------------------------------
# 7. Re-create the hyperlinks that should remain, in the same order
#    they appear in the target file.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:valiantartwear@gamil.com")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:gmail@facebook.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:abc@gmailcom")
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:valiantartwear@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:titovalvala@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A10"), "mailto:titovalvala@gmail.com")

# `Hyperlinks.Add` re-stamps its own "Hyperlink" cellXf on the target
# cell instead of reusing the workbook's existing (identical) one, so
# re-apply the canonical hyperlink-cell format on top to collapse the
# style back onto the one already used by A2/A3/A4's siblings.
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 8. Column width tweaks (D and E got wider).
# ------------------------------------------------------------------
$ws.Range("D1").ColumnWidth = 26.5703125
$ws.Range("E1").ColumnWidth = 16.140625

# ------------------------------------------------------------------
# 9. Selection / active cell moves to B10, and the frozen top-left
#    cell override is dropped.
# ------------------------------------------------------------------
$ws.Range("B1").Select()
$ws.Range("B10").Select()
